$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Martin Guptill"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = " Tabraiz Shamsi"
$ws.Range("J2").Value = "Temba Bavuma(C)"
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = "Bowled"
$ws.Range("N2").Value = " Tim Southee"
$ws.Range("A3").Value = "Daryl Mitchell"
$ws.Range("E3").Value = " Anrich Nortje"
$ws.Range("J3").Value = "Quinton de Kock"
$ws.Range("M3").Value = "Bowled"
$ws.Range("N3").Value = " Trent Boult"
$ws.Range("A4").Value = "Kane Williamson(C)"
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "Caught"
$ws.Range("E4").Value = " Tabraiz Shamsi"
$ws.Range("J4").Value = "Rassie Van der Dussen"
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = "Bowled"
$ws.Range("N4").Value = " Tim Southee"
$ws.Range("A5").Value = "Devon Conway"
$ws.Range("B5").Value = 16
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = "Caught"
$ws.Range("E5").Value = " Kagiso Rabada"
$ws.Range("J5").Value = "Aiden Markram"
$ws.Range("K5").Value = 39
$ws.Range("L5").Value = 15
$ws.Range("M5").Value = "LBW"
$ws.Range("N5").Value = " Tim Southee"
$ws.Range("A6").Value = "Glenn Phillips"
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = "LBW"
$ws.Range("E6").Value = " Dwaine Pretorius"
$ws.Range("J6").Value = "David Miller"
$ws.Range("K6").Value = 36
$ws.Range("L6").Value = 13
$ws.Range("N6").Value = " Trent Boult"
$ws.Range("A7").Value = "James Neesham"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 9
$ws.Range("E7").Value = " Tabraiz Shamsi"
$ws.Range("J7").Value = "Reeza Hendricks"
$ws.Range("K7").Value = 5
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = "Caught"
$ws.Range("N7").Value = " Trent Boult"
$ws.Range("A8").Value = "Mitchell Santner"
$ws.Range("B8").Value = 32
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "LBW"
$ws.Range("E8").Value = " Anrich Nortje"
$ws.Range("J8").Value = "Dwaine Pretorius"
$ws.Range("K8").Value = 3
$ws.Range("N8").Value = " Ish Sodhi"
$ws.Range("A9").Value = "Adam Milne"
$ws.Range("B9").Value = 27
$ws.Range("D9").Value = "Bowled"
$ws.Range("E9").Value = " Kagiso Rabada"
$ws.Range("J9").Value = "Kagiso Rabada"
$ws.Range("K9").Value = 24
$ws.Range("L9").Value = 11
$ws.Range("M9").Value = "Bowled"
$ws.Range("N9").Value = " Adam Milne"
$ws.Range("A10").Value = "Ish Sodhi"
$ws.Range("B10").Value = 34
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = "Caught"
$ws.Range("E10").Value = " Anrich Nortje"
$ws.Range("J10").Value = "Keshav Maharaj"
$ws.Range("K10").Value = 18
$ws.Range("L10").Value = 8
$ws.Range("M10").Value = "NOT OUT"
$ws.Range("N10").Value = " "
$ws.Range("A11").Value = "Tim Southee"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 1
$ws.Range("J11").Value = "Anrich Nortje"
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = "LBW"
$ws.Range("N11").Value = " Trent Boult"
$ws.Range("A12").Value = "Trent Boult"
$ws.Range("B12").Value = 26
$ws.Range("C12").Value = 8
$ws.Range("E12").Value = " Keshav Maharaj"
$ws.Range("J12").Value = "Tabraiz Shamsi"
$ws.Range("K12").Value = 11
$ws.Range("L12").Value = 7
$ws.Range("N12").Value = " Ish Sodhi"
$ws.Range("A16").Value = 191
$ws.Range("C16").Value = "'12.4"
$ws.Range("D16").Value = 76
$ws.Range("J16").Value = 154
$ws.Range("L16").Value = "'12.1"
$ws.Range("M16").Value = 73
$ws.Range("A21").Value = "Kagiso Rabada"
$ws.Range("B21").Value = "'2.0"
$ws.Range("C21").Value = 29
$ws.Range("E21").Value = 14.5
$ws.Range("J21").Value = "Adam Milne"
$ws.Range("K21").Value = "'2.0"
$ws.Range("L21").Value = 33
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 16.5
$ws.Range("A22").Value = "Dwaine Pretorius"
$ws.Range("B22").Value = "'2.0"
$ws.Range("C22").Value = 20
$ws.Range("E22").Value = 10
$ws.Range("J22").Value = "Mitchell Santner"
$ws.Range("K22").Value = "'2.0"
$ws.Range("L22").Value = 15
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 7.5
$ws.Range("A23").Value = "Tabraiz Shamsi"
$ws.Range("D23").Value = 3
$ws.Range("J23").Value = "Trent Boult"
$ws.Range("K23").Value = "'3.0"
$ws.Range("L23").Value = 34
$ws.Range("M23").Value = 4
$ws.Range("N23").Value = 11.33
$ws.Range("A24").Value = "Anrich Nortje"
$ws.Range("C24").Value = 48
$ws.Range("E24").Value = 16
$ws.Range("J24").Value = "Tim Southee"
$ws.Range("K24").Value = "'3.0"
$ws.Range("L24").Value = 40
$ws.Range("M24").Value = 3
$ws.Range("N24").Value = 13.33
$ws.Range("A25").Value = "Keshav Maharaj"
$ws.Range("B25").Value = "'2.4"
$ws.Range("C25").Value = 48
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 20
$ws.Range("J25").Value = "Ish Sodhi"
$ws.Range("K25").Value = "'2.1"
$ws.Range("L25").Value = 32
$ws.Range("N25").Value = 15.24
